# Uploaded data for evaluation apps - participant 1
# Adds three new worksheets ("token", "droidweight", "growtracker"), each
# holding the first 5 trace rows (Execution number / log) copied from the
# existing "GNUCASH-1.0.3" sheet, and updates the active-sheet/selection
# state to match.

$wb = $excel.ActiveWorkbook

# Template sheet whose header + style formatting (bold header row style,
# data-row style) the three new sheets should inherit exactly.
$template = $wb.Worksheets.Item("GNUCASH-1.0.3")

function Add-TraceSheet {
    # Positional params (named `-Param value` binding is unreliable here).
    param($SheetName, $Col1Width, $Col2Width)

    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $afterSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $SheetName

    # Template has 10 data rows (1..10); keep only the first 5 (rows 2..6),
    # drop rows 7..11.
    $newSheet.Rows("7:11").Delete()

    $newSheet.Columns.Item(1).ColumnWidth = $Col1Width
    $newSheet.Columns.Item(2).ColumnWidth = $Col2Width

    $newSheet.Range("A1:B6").Select()

    return $newSheet
}

$tokenSheet = Add-TraceSheet "token" 14.166666666666666 23.998697916666668
$droidweightSheet = Add-TraceSheet "droidweight" 19.498697916666668 20.666666666666668
$growtrackerSheet = Add-TraceSheet "growtracker" 20.998697916666668 20.166666666666668

# Fix up the template sheet's own selection (now highlighting the data
# range it was copied from).
$template.Range("A1:B6").Select()

# "growtracker" (the last-added sheet) ends up the active tab; give it a
# single-cell selection distinct from the A1:B6 block used on the others.
$growtrackerSheet.Activate()
$growtrackerSheet.Range("C5").Select()

Write-Output "done"
